# Auto-generated: updates Moogle Profits market-price / profit figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 14).ClearContents()

# Row 40
$ws.Cells.Item(40, 8).Value = 2192.92
$ws.Cells.Item(40, 9).Value = 1526.2
$ws.Cells.Item(40, 10).Value = 3193
$ws.Cells.Item(40, 11).Value = 1526.2
$ws.Cells.Item(40, 12).Value = 3193
$ws.Cells.Item(40, 13).Value = -1351.2
$ws.Cells.Item(40, 14).Value = -3543

# Row 80
$ws.Cells.Item(80, 8).Value = 1872.1765
$ws.Cells.Item(80, 9).Value = 268.75
$ws.Cells.Item(80, 10).Value = 2365.5386
$ws.Cells.Item(80, 11).Value = 806.25
$ws.Cells.Item(80, 12).Value = 7096.6158
$ws.Cells.Item(80, 13).Value = 191.75
$ws.Cells.Item(80, 14).Value = -9092.6158

# Row 83
$ws.Cells.Item(83, 8).Value = 1872.1765
$ws.Cells.Item(83, 9).Value = 268.75
$ws.Cells.Item(83, 10).Value = 2365.5386
$ws.Cells.Item(83, 11).Value = 2418.75
$ws.Cells.Item(83, 12).Value = 21289.8474
$ws.Cells.Item(83, 13).Value = 2573.25
$ws.Cells.Item(83, 14).Value = -31273.8474

# Row 86
$ws.Cells.Item(86, 8).Value = 5056.3335
$ws.Cells.Item(86, 9).Value = 3751
$ws.Cells.Item(86, 10).Value = 5429.2856
$ws.Cells.Item(86, 11).Value = 3751
$ws.Cells.Item(86, 12).Value = 5429.2856
$ws.Cells.Item(86, 13).Value = -2628
$ws.Cells.Item(86, 14).Value = -7675.2856

# Row 89
$ws.Cells.Item(89, 8).Value = 5056.3335
$ws.Cells.Item(89, 9).Value = 3751
$ws.Cells.Item(89, 10).Value = 5429.2856
$ws.Cells.Item(89, 11).Value = 18755
$ws.Cells.Item(89, 12).Value = 27146.428
$ws.Cells.Item(89, 13).Value = -13139
$ws.Cells.Item(89, 14).Value = -38378.428

# Row 129
$ws.Cells.Item(129, 8).Value = 21065.334
$ws.Cells.Item(129, 9).Value = 21065.334
$ws.Cells.Item(129, 11).Value = 63196.00199999999
$ws.Cells.Item(129, 13).Value = -58196.00199999999

# Row 137
$ws.Cells.Item(137, 8).Value = 2055.0212
$ws.Cells.Item(137, 9).Value = 2140.4614
$ws.Cells.Item(137, 10).Value = 1638.5
$ws.Cells.Item(137, 11).Value = 6421.3842
$ws.Cells.Item(137, 12).Value = 4915.5
$ws.Cells.Item(137, 13).Value = -3871.3842
$ws.Cells.Item(137, 14).Value = -10015.5

# Row 138
$ws.Cells.Item(138, 8).Value = 7372.6523
$ws.Cells.Item(138, 10).Value = 9492.700000000001
$ws.Cells.Item(138, 12).Value = 28478.1
$ws.Cells.Item(138, 14).Value = -38758.10000000001


$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 3284.5715
$ws.Cells.Item(45, 9).Value = 1665.2222
$ws.Cells.Item(45, 11).Value = 1665.2222
$ws.Cells.Item(45, 13).Value = -1288.2222

# Row 61
$ws.Cells.Item(61, 8).Value = 5078.3335
$ws.Cells.Item(61, 9).Value = 4794.6875
$ws.Cells.Item(61, 10).Value = 7347.5
$ws.Cells.Item(61, 11).Value = 4794.6875
$ws.Cells.Item(61, 12).Value = 7347.5
$ws.Cells.Item(61, 13).Value = -4582.6875
$ws.Cells.Item(61, 14).Value = -7771.5

# Row 74
$ws.Cells.Item(74, 8).Value = 8775192
$ws.Cells.Item(74, 9).Value = 10419204
$ws.Cells.Item(74, 11).Value = 10419204
$ws.Cells.Item(74, 13).Value = -10418330

# Row 77
$ws.Cells.Item(77, 8).Value = 8775192
$ws.Cells.Item(77, 9).Value = 10419204
$ws.Cells.Item(77, 11).Value = 52096020
$ws.Cells.Item(77, 13).Value = -52091652

# Row 132
$ws.Cells.Item(132, 8).Value = 3673.0588
$ws.Cells.Item(132, 9).Value = 2150.024
$ws.Cells.Item(132, 11).Value = 6450.072
$ws.Cells.Item(132, 13).Value = -3920.072

# Row 136
$ws.Cells.Item(136, 8).Value = 5078.3335
$ws.Cells.Item(136, 9).Value = 4794.6875
$ws.Cells.Item(136, 10).Value = 7347.5
$ws.Cells.Item(136, 11).Value = 14384.0625
$ws.Cells.Item(136, 12).Value = 22042.5
$ws.Cells.Item(136, 13).Value = -11834.0625
$ws.Cells.Item(136, 14).Value = -27142.5


$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Cells.Item(82, 8).Value = 77843.625
$ws.Cells.Item(82, 9).Value = 5689
$ws.Cells.Item(82, 11).Value = 5689
$ws.Cells.Item(82, 13).Value = -5306

# Row 85
$ws.Cells.Item(85, 8).Value = 77843.625
$ws.Cells.Item(85, 9).Value = 5689
$ws.Cells.Item(85, 11).Value = 5689
$ws.Cells.Item(85, 13).Value = -4363

# Row 99
$ws.Cells.Item(99, 8).Value = 1431.8125
$ws.Cells.Item(99, 9).Value = 1260.6
$ws.Cells.Item(99, 11).Value = 1260.6
$ws.Cells.Item(99, 13).Value = 237.4000000000001

# Row 107
$ws.Cells.Item(107, 8).Value = 6000
$ws.Cells.Item(107, 9).Value = 6000
$ws.Cells.Item(107, 11).Value = 6000
$ws.Cells.Item(107, 13).Value = -4080

# Row 132
$ws.Cells.Item(132, 8).Value = 121079.836
$ws.Cells.Item(132, 10).Value = 121079.836
$ws.Cells.Item(132, 12).Value = 121079.836
$ws.Cells.Item(132, 14).Value = -131199.836

# Row 134
$ws.Cells.Item(134, 8).Value = 1893.8
$ws.Cells.Item(134, 9).Value = 1420.0883
$ws.Cells.Item(134, 11).Value = 4260.2649
$ws.Cells.Item(134, 13).Value = -1725.2649


$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 2340.3333
$ws.Cells.Item(16, 9).Value = 2340.3333
$ws.Cells.Item(16, 11).Value = 2340.3333
$ws.Cells.Item(16, 13).Value = -2053.3333

# Row 31
$ws.Cells.Item(31, 8).Value = 7628.1
$ws.Cells.Item(31, 10).Value = 15649.091
$ws.Cells.Item(31, 12).Value = 15649.091
$ws.Cells.Item(31, 14).Value = -16239.091

# Row 34
$ws.Cells.Item(34, 8).Value = 7628.1
$ws.Cells.Item(34, 10).Value = 15649.091
$ws.Cells.Item(34, 12).Value = 15649.091
$ws.Cells.Item(34, 14).Value = -16053.091

# Row 37
$ws.Cells.Item(37, 8).Value = 5000
$ws.Cells.Item(37, 9).Value = 5000
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 5000
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 13).Value = -4893
$ws.Cells.Item(37, 14).ClearContents()

# Row 38
$ws.Cells.Item(38, 8).Value = 17124.75
$ws.Cells.Item(38, 9).Value = 9500
$ws.Cells.Item(38, 10).Value = 39999
$ws.Cells.Item(38, 11).Value = 9500
$ws.Cells.Item(38, 12).Value = 39999
$ws.Cells.Item(38, 13).Value = -9123
$ws.Cells.Item(38, 14).Value = -40753

# Row 41
$ws.Cells.Item(41, 8).Value = 59999
$ws.Cells.Item(41, 10).Value = 59999
$ws.Cells.Item(41, 12).Value = 59999
$ws.Cells.Item(41, 14).Value = -60855

# Row 46
$ws.Cells.Item(46, 8).Value = 17124.75
$ws.Cells.Item(46, 9).Value = 9500
$ws.Cells.Item(46, 10).Value = 39999
$ws.Cells.Item(46, 11).Value = 9500
$ws.Cells.Item(46, 12).Value = 39999
$ws.Cells.Item(46, 13).Value = -9289
$ws.Cells.Item(46, 14).Value = -40421

# Row 55
$ws.Cells.Item(55, 8).Value = 33666.668
$ws.Cells.Item(55, 9).Value = 25000
$ws.Cells.Item(55, 11).Value = 25000
$ws.Cells.Item(55, 13).Value = -24685

# Row 99
$ws.Cells.Item(99, 8).Value = 1394.1052
$ws.Cells.Item(99, 9).Value = 1363.0605
$ws.Cells.Item(99, 11).Value = 1363.0605
$ws.Cells.Item(99, 13).Value = 134.9395

# Row 113
$ws.Cells.Item(113, 8).Value = 2340.3333
$ws.Cells.Item(113, 9).Value = 2340.3333
$ws.Cells.Item(113, 11).Value = 2340.3333
$ws.Cells.Item(113, 13).Value = -170.3332999999998

# Row 122
$ws.Cells.Item(122, 8).Value = 1618.8462
$ws.Cells.Item(122, 9).Value = 1670.6666
$ws.Cells.Item(122, 10).Value = 1502.25
$ws.Cells.Item(122, 11).Value = 5011.9998
$ws.Cells.Item(122, 12).Value = 4506.75
$ws.Cells.Item(122, 13).Value = -2561.9998
$ws.Cells.Item(122, 14).Value = -9406.75

# Row 126
$ws.Cells.Item(126, 8).Value = 1394.1052
$ws.Cells.Item(126, 9).Value = 1363.0605
$ws.Cells.Item(126, 11).Value = 4089.1815
$ws.Cells.Item(126, 13).Value = -1619.1815


$ws = $wb.Worksheets.Item("CUL")
# Row 112
$ws.Cells.Item(112, 8).Value = 13127.091
$ws.Cells.Item(112, 9).Value = 7566.5
$ws.Cells.Item(112, 11).Value = 22699.5
$ws.Cells.Item(112, 13).Value = -21591.5

# Row 116
$ws.Cells.Item(116, 8).Value = 41130
$ws.Cells.Item(116, 9).Value = 1483.3334
$ws.Cells.Item(116, 10).Value = 100600
$ws.Cells.Item(116, 11).Value = 4450.0002
$ws.Cells.Item(116, 12).Value = 301800
$ws.Cells.Item(116, 13).Value = -1008.0002
$ws.Cells.Item(116, 14).Value = -308684

# Row 137
$ws.Cells.Item(137, 8).Value = 3519.963
$ws.Cells.Item(137, 10).Value = 4053.1667
$ws.Cells.Item(137, 12).Value = 12159.5001
$ws.Cells.Item(137, 14).Value = -22359.5001


$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 11582.333
$ws.Cells.Item(80, 9).Value = 6499.6665
$ws.Cells.Item(80, 11).Value = 6499.6665
$ws.Cells.Item(80, 13).Value = -5501.6665

# Row 83
$ws.Cells.Item(83, 8).Value = 11582.333
$ws.Cells.Item(83, 9).Value = 6499.6665
$ws.Cells.Item(83, 11).Value = 32498.3325
$ws.Cells.Item(83, 13).Value = -27506.3325


$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 7561.0713
$ws.Cells.Item(40, 9).Value = 6612.9585
$ws.Cells.Item(40, 11).Value = 6612.9585
$ws.Cells.Item(40, 13).Value = -6476.9585

# Row 43
$ws.Cells.Item(43, 8).Value = 209996.67
$ws.Cells.Item(43, 10).Value = 209996.67
$ws.Cells.Item(43, 12).Value = 209996.67
$ws.Cells.Item(43, 14).Value = -210382.67

# Row 132
$ws.Cells.Item(132, 8).Value = 3423
$ws.Cells.Item(132, 9).Value = 999.8570999999999
$ws.Cells.Item(132, 11).Value = 2999.5713
$ws.Cells.Item(132, 13).Value = -469.5712999999996


$ws = $wb.Worksheets.Item("WVR")
# Row 33
$ws.Cells.Item(33, 8).Value = 14500
$ws.Cells.Item(33, 9).Value = 19000
$ws.Cells.Item(33, 11).Value = 19000
$ws.Cells.Item(33, 13).Value = -18750

# Row 36
$ws.Cells.Item(36, 8).Value = 14500
$ws.Cells.Item(36, 9).Value = 19000
$ws.Cells.Item(36, 11).Value = 19000
$ws.Cells.Item(36, 13).Value = -18750

# Row 113
$ws.Cells.Item(113, 8).Value = 975.5
$ws.Cells.Item(113, 9).Value = 846.8182
$ws.Cells.Item(113, 11).Value = 2540.4546
$ws.Cells.Item(113, 13).Value = -370.4546

# Row 126
$ws.Cells.Item(126, 8).Value = 1933.6394
$ws.Cells.Item(126, 9).Value = 1819.6735
$ws.Cells.Item(126, 10).Value = 2399
$ws.Cells.Item(126, 11).Value = 5459.020500000001
$ws.Cells.Item(126, 12).Value = 7197
$ws.Cells.Item(126, 13).Value = -2989.020500000001
$ws.Cells.Item(126, 14).Value = -12137

